$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relabel the four "Layer" quadrant headers (each shifts down by one: 1->0, 2->1, 3->2, 4->3) ---
$ws.Range("A1").Value = "Layer 0"
$ws.Range("G1").Value = "Layer 1"
$ws.Range("M1").Value = "Layer 2"
$ws.Range("S1").Value = "Layer 3"

# --- Correct the "Layer 1" (columns G:K) print-time data, rows 3-14 ---
$ws.Range("J3").Value = 19

$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 17

$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 38

$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 18
$ws.Range("J6").Value = 8

$ws.Range("H7").Value = 2
$ws.Range("I7").Value = 7
$ws.Range("J7").Value = 6

$ws.Range("I8").Value = 4
$ws.Range("J8").Value = 26

$ws.Range("I9").Value = 10
$ws.Range("J9").Value = 39

$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 22
$ws.Range("J10").Value = 3

$ws.Range("I11").Value = 2
$ws.Range("J11").Value = 54

$ws.Range("I12").Value = 8
$ws.Range("J12").Value = 48

$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 10
$ws.Range("J13").Value = 49

$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 29

# --- Update the active cell selection to match the author's final cursor position ---
$ws.Range("H15").Select()
